$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.976.26'
$ws.Range('E2').Value = '  -0.77%  '
$ws.Range('D3').Value = '1.819.42'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').Value = "'311.39"
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').Value = "'0.4526"
$ws.Range('E7').Value = '  +6.43%  '
$ws.Range('D8').Value = "'0.3698"
$ws.Range('E8').Value = '  +0.64%  '
$ws.Range('D9').Value = "'0.07282"
$ws.Range('E9').Value = '  +0.81%  '
$ws.Range('D10').Value = "'0.8547"
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('D11').Value = "'20.69"
$ws.Range('E11').Value = '  -1.21%  '
$ws.Range('D12').Value = '1.807.12'
$ws.Range('E12').Value = '  -0.74%  '
$ws.Range('D13').Value = "'6.643"
$ws.Range('E13').Value = '  -0.34%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'5.328"
$ws.Range('E14').Value = '  +0.49%  '
$ws.Range('D15').Value = "'0.07094"
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = "'92.08"
$ws.Range('E16').Value = '  +4.63%  '
$ws.Range('D17').Value = "'1.002"
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('D18').Value = "'0.000008784"
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('D20').Value = "'14.96"
$ws.Range('E20').Value = '  -0.60%  '
$ws.Range('D21').Value = '26.984.97'
$ws.Range('E21').Value = '  -0.89%  '
$ws.Range('D22').Value = "'5.161"
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('D23').Value = "'10.94"
$ws.Range('E23').Value = '  +0.92%  '
$ws.Range('D24').Value = "'1.986"
$ws.Range('E24').Value = '  -0.85%  '
$ws.Range('D25').Value = "'151.64"
$ws.Range('E25').Value = '  -1.21%  '
$ws.Range('D26').Value = "'2.210"
$ws.Range('E26').Value = '  +4.63%  '
$ws.Range('D27').Value = "'18.45"
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('D28').Value = "'5.235"
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('D29').Value = "'116.40"
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('D30').Value = "'0.08869"
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').Value = "'1.181"
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').Value = "'0.7518"
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('D33').Value = "'2.965"
$ws.Range('E33').Value = '  +5.65%  '
$ws.Range('D34').Value = "'4.430"
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('E35').Value = '  -0.49%  '
$ws.Range('D36').Value = "'1.100"
$ws.Range('E36').Value = '  -0.95%  '
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').Value = "'0.05237"
$ws.Range('E38').Value = '  -0.50%  '
$ws.Range('D39').Value = "'0.5291"
$ws.Range('E39').Value = '  +5.16%  '
$ws.Range('D40').Value = "'7.170"
$ws.Range('E40').Value = '  +0.70%  '
$ws.Range('D41').Value = "'2.872"
$ws.Range('E41').Value = '  +0.32%  '
$ws.Range('D42').Value = "'0.1707"
$ws.Range('E42').Value = '  +1.19%  '
$ws.Range('D43').Value = "'0.5195"
$ws.Range('E43').Value = '  +9.78%  '
$ws.Range('D44').Value = "'8.504"
$ws.Range('E44').Value = '  -1.08%  '
$ws.Range('D45').Value = "'10.56"
$ws.Range('E45').Value = '  -0.58%  '
$ws.Range('D46').Value = "'1.967"
$ws.Range('E46').Value = '  +9.20%  '
$ws.Range('D47').Value = "'105.62"
$ws.Range('E47').Value = '  -1.04%  '
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('D49').Value = "'1.665"
$ws.Range('E49').Value = '  +0.62%  '
$ws.Range('D50').Value = "'0.06371"
$ws.Range('E50').Value = '  +0.07%  '
$ws.Range('D51').Value = "'0.9165"
$ws.Range('E51').Value = '  +0.32%  '
